$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new timesheet entry row 7
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7").Value = 43857
$ws.Range("B7").Value = 3.5
$ws.Range("C7").Value = "7pm"
$ws.Range("D7").Value = "10:30pm"

# Update selection to match the new active cell after editing
$ws.Range("C14").Select()

$wb.Save()
